$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 6262, 46078.95833333334),
    @(3, 6217, 46078.96875),
    @(4, 6168, 46078.97916666666),
    @(5, 6113, 46078.98958333334),
    @(6, 6125, 46079),
    @(7, 6048, 46079.01041666666),
    @(8, 6031, 46079.02083333334),
    @(9, 5978, 46079.04166666666),
    @(10, 5929, 46079.05208333334),
    @(11, 5912, 46079.0625),
    @(12, 5898, 46079.07291666666),
    @(13, 5886, 46079.08333333334),
    @(14, 5895, 46079.09375),
    @(15, 5856, 46079.10416666666),
    @(16, 5925, 46079.11458333334),
    @(17, 5924, 46079.125),
    @(18, 5955, 46079.13541666666),
    @(19, 5935, 46079.14583333334),
    @(20, 5960, 46079.15625),
    @(21, 6087, 46079.16666666666),
    @(22, 6137, 46079.17708333334),
    @(23, 6226, 46079.1875),
    @(24, 6328, 46079.19791666666),
    @(25, 6546, 46079.20833333334),
    @(26, 6707, 46079.21875),
    @(27, 6791, 46079.22916666666),
    @(28, 6971, 46079.23958333334),
    @(29, 7154, 46079.25),
    @(30, 7393, 46079.26041666666),
    @(31, 7498, 46079.27083333334),
    @(32, 7558, 46079.28125),
    @(33, 7633, 46079.29166666666),
    @(34, 7648, 46079.30208333334),
    @(35, 7607, 46079.3125),
    @(36, 7569, 46079.32291666666),
    @(37, 7562, 46079.33333333334),
    @(38, 7518, 46079.34375),
    @(39, 7459, 46079.35416666666),
    @(40, 7365, 46079.36458333334),
    @(41, 7149, 46079.375),
    @(42, 7015, 46079.38541666666),
    @(43, 6910, 46079.39583333334)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value2 = $item[1]
    $ws.Cells.Item($r, 2).Value2 = $item[2]
}

# Ensure the timestamp column keeps its custom date/time number format
# for the newly added rows (35-43), matching the existing rows above.
$fmt = $ws.Range("B34").NumberFormat
$ws.Range("B35:B43").NumberFormat = $fmt
